$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text-valued (shared-string) cells: Repayment_amount and Pending Amount Recovery ---
# These cells hold text that looks numeric ("1,044,908.00", "0.68", ...). Assigning a
# plain string would let Excel auto-convert it into a real number, so we force the
# cell to Text format first, write the string, then restore a plain/unformatted style
# so the cell's number format matches the rest of the sheet (General).

function Set-TextValue($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

# Erick Ervan Dewanggga (row 2)
Set-TextValue "E2" "1,044,908.00"
Set-TextValue "G2" "0.68"

# Sucika Wardani (row 3)
Set-TextValue "E3" "1,675,977.00"
Set-TextValue "G3" "1.08"

# Axl Wicaksono (row 6)
Set-TextValue "E6" "1,218,637.00"
Set-TextValue "G6" "0.90"

# Riska Nurlita (row 9)
Set-TextValue "E9" "834,007.00"
Set-TextValue "G9" "0.44"

# Yandi Nugraha (row 16)
Set-TextValue "E16" "2,234,658.00"
Set-TextValue "G16" "1.74"

# --- Numeric Repayment_collections counts (column D) ---
$ws.Range("D2").Value = 5
$ws.Range("D3").Value = 6
$ws.Range("D6").Value = 4
$ws.Range("D9").Value = 6
$ws.Range("D16").Value = 9

# --- Rename the sheet (revision bump: (5) -> (6)) ---
$ws.Name = "repayment_20250915_20250915 (6)"
